$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared "MagicBallController" string to "MagicBall" - this is
# the InitWeapon value used by BoneMan/Reaper/Balance (column L, rows 4-6).
$ws.Cells.Replace("MagicBallController", "MagicBall")

# Left-align the AttackCooldown values (J4:J6) so they match the rest of
# the data row formatting (same look as the PrefabPath column).
$ws.Range("J4:J6").HorizontalAlignment = -4131

# Give the InitWeapon cells (L4:L6) their own distinct style (reading
# order explicitly set, like the rest of the sheet).
$ws.Range("L4:L6").ReadingOrder = 1
$ws.Range("L4:L6").ShrinkToFit = $true
